$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.650.17'
$ws.Range("E2").Value = '  +2.12%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.868.73'
$ws.Range("E3").Value = '  +1.14%  '

$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.83'
$ws.Range("E5").Value = '  +1.15%  '

$ws.Range("E6").Value = '  -0.05%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4825'
$ws.Range("E7").Value = '  +1.25%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3811'
$ws.Range("E8").Value = '  +3.48%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07368'
$ws.Range("E9").Value = '  +1.86%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9375'
$ws.Range("E10").Value = '  +0.59%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.90'
$ws.Range("E11").Value = '  +5.09%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07804'
$ws.Range("E12").Value = '  +0.27%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.889.81'
$ws.Range("E13").Value = '  +2.40%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.468'
$ws.Range("E14").Value = '  +1.45%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.580'
$ws.Range("E15").Value = '  +1.53%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '90.41'
$ws.Range("E16").Value = '  +1.42%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008854'
$ws.Range("E18").Value = '  +1.92%  '

$ws.Range("E19").Value = '  -0.11%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.724.66'
$ws.Range("E20").Value = '  +2.27%  '

$ws.Range("E21").Value = '  +1.03%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.115'
$ws.Range("E22").Value = '  +1.31%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.108.61'
$ws.Range("E23").Value = '  +1.44%  '

$ws.Range("E24").Value = '  +1.16%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.936'
$ws.Range("E25").Value = '  +0.12%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.53'
$ws.Range("E26").Value = '  +2.34%  '

$ws.Range("E27").Value = '  +1.07%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.059'
$ws.Range("E28").Value = '  +3.73%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '115.72'
$ws.Range("E29").Value = '  +1.00%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.959'
$ws.Range("E30").Value = '  +1.34%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08918'
$ws.Range("E31").Value = '  +0.67%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.338'
$ws.Range("E32").Value = '  +0.69%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.214'
$ws.Range("E33").Value = '  +3.08%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7611'
$ws.Range("E34").Value = '  +3.14%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.618'
$ws.Range("E35").Value = '  +2.25%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.709'
$ws.Range("E36").Value = '  +0.60%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.137'
$ws.Range("E37").Value = '  +2.19%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02041'
$ws.Range("E38").Value = '  +3.34%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5644'
$ws.Range("E39").Value = '  +7.38%  '

$ws.Range("E40").Value = '  +2.35%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.985'
$ws.Range("E41").Value = '  +0.58%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.066'
$ws.Range("E42").Value = '  +0.64%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.558'
$ws.Range("E43").Value = '  +3.37%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1530'
$ws.Range("E44").Value = '  +0.63%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.78'
$ws.Range("E45").Value = '  +2.14%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4905'
$ws.Range("E46").Value = '  +3.65%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.013'
$ws.Range("E47").Value = '  -0.08%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '104.96'
$ws.Range("E48").Value = '  +2.87%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.673'
$ws.Range("E49").Value = '  +3.62%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '67.73'
$ws.Range("E50").Value = '  +2.76%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06102'
$ws.Range("E51").Value = '  +0.84%  '
